# This revision's unified diff only touches two SharePoint-managed
# customXml parts:
#
#   - customXml/item2.xml   (ct:contentTypeSchema) — its
#     ma:versionID and ma:fieldsID attributes are replaced with new
#     server-minted hex tokens.
#   - customXml/itemProps2.xml — its ds:itemID GUID is replaced with a
#     freshly minted GUID and a ds:schemaRefs block is (re)materialized.
#
# These values are random identifiers stamped by the SharePoint content
# type hub when a library's content type is republished/synced; they
# are not derived from, or tied to, any visible document content, and
# Word does not expose a UI/automation surface that lets a user (or a
# macro) set them directly — CustomXMLParts in the object model is
# read-only in this respect (there is no supported way to assign new
# GUIDs/version stamps to the SharePoint content-type schema part).
#
# Consistent with that: there is no matching content edit anywhere in
# this document either — no hyperlink, field, or visible "OP10" text
# exists in word/document.xml, headers, footers, or the glossary part
# that a Find/Replace could target. The commit's "Fixed broken link in
# OP10" message describes a fix made to a different document (OP10) in
# the same repository/commit; this file simply picked up the
# SharePoint content-type resync as a side effect of being saved
# alongside it.
#
# So the faithful reproduction here is simply to resave the document
# as-is — no body/header/footer/style edit is warranted because none
# occurred.
$d = $word.ActiveDocument
$d.Save()
